# The edit re-orders the five observation records currently stored in
# rows 2-6 of the active sheet. Every column is identical between the
# "old" and "new" position of a record except for:
#   A  (Id), B (Taxonsorteringsordning), D (Rödlistade), E (TaxonId),
#   F  (Artnamn), G (Vetenskapligt namn), H (Auktor),
#   Q  (Ost), R (Nord), Z (Starttid), AB (Sluttid)
# so only those columns actually need to move between rows.
#
# New row  <- old row
#   2      <- 6
#   3      <- 5
#   4      <- 2
#   5      <- 3
#   6      <- 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")
$firstRow = 2
$lastRow = 6

# Snapshot the current ("before") values for every relevant cell so that
# the subsequent writes don't clobber data we still need to read.
$snapshot = @{}
foreach ($col in $cols) {
    $colValues = @{}
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $colValues[$r] = $ws.Range("$col$r").Value2
    }
    $snapshot[$col] = $colValues
}

# New row number -> source (old) row number.
$rowMap = @{ 2 = 6; 3 = 5; 4 = 2; 5 = 3; 6 = 4 }

foreach ($col in $cols) {
    foreach ($newRow in $rowMap.Keys) {
        $oldRow = $rowMap[$newRow]
        $ws.Range("$col$newRow").Value2 = $snapshot[$col][$oldRow]
    }
}
